# The deck ships with the "Integral" (Red Violet) theme applied to the
# slide master/design (ppt/theme/theme2.xml, the theme actually wired up
# via presentation.xml.rels + slideMaster1.xml.rels). The edit swaps the
# applied theme's colour scheme back to the stock Office palette (the
# colours that used to live, unused, in ppt/theme/theme1.xml), i.e. the
# equivalent of Design > Variants > Colors > "Office" in the UI.
#
# RGB() values below are plain VBA-style BGR-packed integers
# (R + G*256 + B*65536) for each of the 12 standard theme colour slots,
# in PowerPoint's fixed ColorScheme.Colors(index) order:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2  7 accent3
#   8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink

$p = $ppt.ActivePresentation

$officeColors = @(
    0,          # 1  dk1      000000
    16777215,   # 2  lt1      FFFFFF
    6968388,    # 3  dk2      44546A
    15132391,   # 4  lt2      E7E6E6
    13998939,   # 5  accent1  5B9BD5
    3243501,    # 6  accent2  ED7D31
    10855845,   # 7  accent3  A5A5A5
    49407,      # 8  accent4  FFC000
    12874308,   # 9  accent5  4472C4
    4697456,    # 10 accent6  70AD47
    12673797,   # 11 hlink    0563C1
    7491477     # 12 folHlink 954F72
)

$colorScheme = $p.SlideMaster.ColorScheme
for ($i = 0; $i -lt $officeColors.Length; $i++) {
    $colorScheme.Colors($i + 1).RGB = $officeColors[$i]
}
